# Apply Coinranking crypto snapshot update (GitHub Actions refresh)
# Each touched cell is forced to Text format before the write so the
# literal display string (thousand-dot prices, padded percentages,
# subscript-zero prices, etc.) survives instead of being re-interpreted
# as a number/date by Excel's automatic type inference.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "35.060.75"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.36%  "

# Row 3: Ethereum
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.856.02"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +2.64%  "

# Row 4: TetherUSD
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.06%  "

# Row 5: BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.64"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +2.39%  "

# Row 6: XRP
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.623"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +2.11%  "

# Row 7: USDC
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.06%  "

# Row 8: Solana
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "42.82"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +9.38%  "

# Row 9: Cardano
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +2.29%  "

# Row 10: Dogecoin
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0695"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +2.24%  "

# Row 11: TRON
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0988"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.28%  "

# Row 12: WrappedliquidstakedEther2.0
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.128.20"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +2.86%  "

# Row 13: WrappedEther -> Chainlink
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.42"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +2.22%  "

# Row 14: Chainlink -> WrappedEther
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.855.63"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +2.27%  "

# Row 15: Polygon
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.681"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +2.84%  "

# Row 16: Polkadot
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +2.65%  "

# Row 17: WrappedBTC
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "35.042.50"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.26%  "

# Row 18: Litecoin
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "70.36"

# Row 19: ShibaInu
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0796"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.86%  "

# Row 20: BitcoinCash
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "241.20"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.06%  "

# Row 21: Avalanche
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.21"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +2.93%  "

# Row 22: Uniswap
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.76"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +2.34%  "

# Row 23: Dai
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.08%  "

# Row 24: Toncoin
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.30"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +2.76%  "

# Row 25: Monero
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "171.63"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.16%  "

# Row 26: PancakeSwap
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.84"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +24.03%  "

# Row 27: Cosmos
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.92"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +2.08%  "

# Row 28: EthereumClassic
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.68"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +2.82%  "

# Row 29: Stellar
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +2.68%  "

# Row 30: Hedera -> BinanceUSD
$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = "BinanceUSD"
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.01"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.14%  "

# Row 31: BinanceUSD -> Hedera
$ws.Range("B31").NumberFormat = "@"
$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0557"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +2.67%  "

# Row 32: Filecoin
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.00"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.27%  "

# Row 33: InternetComputer(DFINITY)
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.00"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +2.55%  "

# Row 34: LidoDAOToken
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +13.66%  "

# Row 35: WEMIXToken
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +23.23%  "

# Row 36: ImmutableX
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.781"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +12.60%  "

# Row 37: TrustWalletToken
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.22%  "

# Row 38: ARBITRUM
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.08"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +13.07%  "

# Row 39: Aave
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "91.67"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.19%  "

# Row 40: VeChain
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +5.72%  "

# Row 41: Maker
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.349.49"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +2.12%  "

# Row 42: InjectiveProtocol
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "15.08"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +5.77%  "

# Row 43: RenderToken
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +5.79%  "

# Row 44: Gas
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "12.79"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +63.48%  "

# Row 45: HuobiToken
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -2.43%  "

# Row 46: MXToken
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.75"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.13%  "

# Row 47: Kaspa
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0547"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +6.68%  "

# Row 48: FraxShare
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +3.24%  "

# Row 49: RocketPoolETH
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.037.94"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +2.27%  "

# Row 50: Cronos
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +2.60%  "

# Row 51: THORChain
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.43"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +15.69%  "

